$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 146; all existing rows 146:186 shift down to 147:187
$ws.Rows.Item(146).Insert()

# Populate the new row 146 with the new weekly price record
$ws.Cells.Item(146, 1).Value = 4
$ws.Cells.Item(146, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(146, 3).Value = 'Los Lagos'
$ws.Cells.Item(146, 4).Value = 44736
$ws.Cells.Item(146, 5).Value = 10
$ws.Cells.Item(146, 6).Value = 100112009
$ws.Cells.Item(146, 7).Value = 'Acelga'
$ws.Cells.Item(146, 8).Value = 'Sin especificar'
$ws.Cells.Item(146, 9).Value = 'Primera'
$ws.Cells.Item(146, 10).Value = 120
$ws.Cells.Item(146, 11).Value = 12000
$ws.Cells.Item(146, 12).Value = 12000
$ws.Cells.Item(146, 13).Value = 12000
$ws.Cells.Item(146, 14).Value = '$/docena de atados (12 kilos)'
$ws.Cells.Item(146, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(146, 16).Value = 1000
$ws.Cells.Item(146, 17).Value = 12
$ws.Cells.Item(146, 18).Value = 'Hortaliza'
